# Modelo_Licenca_Autorizacion_SIRHA.docx - cosmetic fixes
# 1) Remove a duplicated empty paragraph that was left right after the
#    header table (two identical blank "Normal" paragraphs in a row ->
#    keep only one).
# 2) Bump the left cell margin of the header table from 253 -> 263 twips
#    (12.65pt -> 13.15pt).
# 3) Shrink the row heights of the three "O DIRECTOR GERAL / AVERBAMENTOS"
#    rows (1701 twips -> 1621/1587/1587, i.e. 85.05pt -> 81.05/79.35/79.35pt).
# 4) Register the new ListLabel character styles (316-333) used by the
#    updated list numbering.

$d = $word.ActiveDocument

# --- 1) Drop the duplicate empty paragraph after the header table -----
# (done before touching Tables/Rows so paragraph Range anchors stay valid)
for ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {
    $cur = $d.Paragraphs.Item($i)
    $nxt = $d.Paragraphs.Item($i + 1)
    $curText = $cur.Range.Text
    if ($curText.Length -eq 1 -and $curText -eq $nxt.Range.Text) {
        $cur.Range.Delete()
        break
    }
}

# --- 2) Header table cell margin -------------------------------------
$headerTable = $d.Tables.Item(1)
$headerTable.LeftPadding = 13.15   # 263 twips

# --- 3) Row heights on the AVERBAMENTOS table --------------------------
foreach ($tbl in $d.Tables) {
    if ($tbl.Rows.Count -eq 4 -and $tbl.Range.Text -like "*DIRECTOR GERAL*") {
        $tbl.Rows.Item(2).Height = 81.05   # 1621 twips
        $tbl.Rows.Item(3).Height = 79.35   # 1587 twips
        $tbl.Rows.Item(4).Height = 79.35   # 1587 twips
    }
}

# --- 4) New ListLabel character styles 316-333 -------------------------
for ($i = 316; $i -le 333; $i++) {
    $styleName = "ListLabel$i"
    $style = $d.Styles.Add($styleName, 2)
    $style.NameLocal = "ListLabel $i"
    $style.QuickStyle = $true
    $style.Font.NameBi = "OpenSymbol"
}

Write-Output "done"
